$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.702.67"
$ws.Range("E2").Value = "  +2.98%  "
$ws.Range("D3").Value = "3.208.13"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.206.97"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.32%  "
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E12").Value = "  +4.46%  "
$ws.Range("E13").Value = "  +3.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.11%  "
$ws.Range("D15").Value = "3.732.05"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "66.583.39"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("E17").Value = "  +5.31%  "
$ws.Range("D18").Value = "3.206.29"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "514.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.741"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.53%  "
$ws.Range("E23").Value = "  +5.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  +5.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +16.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("E33").Value = "  +4.24%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "491.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0898"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0423"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.123"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.22%  "
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("E43").Value = "  +6.68%  "
$ws.Range("D44").Value = "0.0₃0671"
$ws.Range("E44").Value = "  +17.22%  "
$ws.Range("D45").Value = "2.935.31"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("E50").Value = "  +4.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.35%  "
